# Adição de correção deriva instrumental
# Apply incremental "instrumental drift" correction to column L (Horas) for rows 4-12,
# and update the active selection to L16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value = 9
$ws.Range("L5").Value = 10
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 12
$ws.Range("L8").Value = 13
$ws.Range("L9").Value = 14
$ws.Range("L10").Value = 15
$ws.Range("L11").Value = 16
$ws.Range("L12").Value = 17

$ws.Range("L16").Select()
